$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3744.5
$ws.Range("I38").Value = 999
$ws.Range("J38").Value = 6490
$ws.Range("K38").Value = 2997
$ws.Range("L38").Value = 19470
$ws.Range("M38").Value = -2625
$ws.Range("N38").Value = -20214
$ws.Range("H57").Value = 125000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 125000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 375000
$ws.Range("N57").Value = -375998
$ws.Range("M57").ClearContents()
$ws.Range("H58").Value = 6941.4287
$ws.Range("J58").Value = 9383.200000000001
$ws.Range("L58").Value = 28149.6
$ws.Range("N58").Value = -28449.6
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H86").Value = 2900
$ws.Range("I86").Value = 2850
$ws.Range("K86").Value = 2850
$ws.Range("M86").Value = -1727
$ws.Range("H89").Value = 2900
$ws.Range("I89").Value = 2850
$ws.Range("K89").Value = 14250
$ws.Range("M89").Value = -8634
$ws.Range("H132").Value = 2306.75
$ws.Range("I132").Value = 2306.75
$ws.Range("K132").Value = 6920.25
$ws.Range("M132").Value = -4390.25
$ws.Range("H135").Value = 3347.25
$ws.Range("I135").Value = 3347.25
$ws.Range("K135").Value = 30125.25
$ws.Range("M135").Value = -27590.25
$ws.Range("H137").Value = 6020
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 6020
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 18060
$ws.Range("N137").Value = -23160
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3036.625
$ws.Range("I32").Value = 2072.4333
$ws.Range("K32").Value = 2072.4333
$ws.Range("M32").Value = -1785.4333
$ws.Range("H74").Value = 1344.6666
$ws.Range("I74").Value = 1017.25
$ws.Range("K74").Value = 1017.25
$ws.Range("M74").Value = -143.25
$ws.Range("H77").Value = 1344.6666
$ws.Range("I77").Value = 1017.25
$ws.Range("K77").Value = 5086.25
$ws.Range("M77").Value = -718.25
$ws.Range("H97").Value = 234.23077
$ws.Range("I97").Value = 234.23077
$ws.Range("K97").Value = 234.23077
$ws.Range("M97").Value = 261.76923

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 19998
$ws.Range("I20").Value = 19998
$ws.Range("K20").Value = 19998
$ws.Range("M20").Value = -19751

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 7500
$ws.Range("I45").Value = 7500
$ws.Range("K45").Value = 7500
$ws.Range("M45").Value = -6907
$ws.Range("H132").Value = 4648.8887
$ws.Range("I132").Value = 3974
$ws.Range("K132").Value = 11922
$ws.Range("M132").Value = -9392
$ws.Range("H139").Value = 59000
$ws.Range("J139").Value = 59000
$ws.Range("L139").Value = 59000
$ws.Range("N139").Value = -69280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3589.4583
$ws.Range("J122").Value = 3742.1333
$ws.Range("L122").Value = 33679.1997
$ws.Range("N122").Value = -38579.1997
$ws.Range("H129").Value = 1189
$ws.Range("J129").Value = 1880
$ws.Range("L129").Value = 5640
$ws.Range("N129").Value = -15640
$ws.Range("H138").Value = 2400
$ws.Range("I138").Value = 2000
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 6000
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = -860
$ws.Range("N138").Value = -19280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 99.5
$ws.Range("I2").Value = 99.5
$ws.Range("K2").Value = 99.5
$ws.Range("M2").Value = 13.5
$ws.Range("H43").Value = 13055
$ws.Range("I43").Value = 6500
$ws.Range("J43").Value = 14147.5
$ws.Range("K43").Value = 6500
$ws.Range("L43").Value = 14147.5
$ws.Range("M43").Value = -6349
$ws.Range("N43").Value = -14449.5
$ws.Range("H57").Value = 37945
$ws.Range("J57").Value = 37945
$ws.Range("L57").Value = 37945
$ws.Range("N57").Value = -39585
$ws.Range("H80").Value = 3456.8125
$ws.Range("I80").Value = 2030.7
$ws.Range("J80").Value = 5833.6665
$ws.Range("K80").Value = 2030.7
$ws.Range("L80").Value = 5833.6665
$ws.Range("M80").Value = -1032.7
$ws.Range("N80").Value = -7829.6665
$ws.Range("H83").Value = 3456.8125
$ws.Range("I83").Value = 2030.7
$ws.Range("J83").Value = 5833.6665
$ws.Range("K83").Value = 10153.5
$ws.Range("L83").Value = 29168.3325
$ws.Range("M83").Value = -5161.5
$ws.Range("N83").Value = -39152.3325
$ws.Range("H97").Value = 3021.6
$ws.Range("I97").Value = 3021.6
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 3021.6
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -2525.6
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 7501
$ws.Range("I20").Value = 9000
$ws.Range("K20").Value = 9000
$ws.Range("M20").Value = -8774
$ws.Range("H22").Value = 2992.2
$ws.Range("I22").Value = 2832
$ws.Range("J22").Value = 3900
$ws.Range("K22").Value = 2832
$ws.Range("L22").Value = 3900
$ws.Range("M22").Value = -2537
$ws.Range("N22").Value = -4490
$ws.Range("H27").Value = 2992.2
$ws.Range("I27").Value = 2832
$ws.Range("J27").Value = 3900
$ws.Range("K27").Value = 2832
$ws.Range("L27").Value = 3900
$ws.Range("M27").Value = -2725
$ws.Range("N27").Value = -4114
$ws.Range("H61").Value = 2841
$ws.Range("I61").Value = 2481.5
$ws.Range("K61").Value = 2481.5
$ws.Range("M61").Value = -2279.5
$ws.Range("H82").Value = 1100.5
$ws.Range("J82").Value = 134
$ws.Range("L82").Value = 134
$ws.Range("N82").Value = -856
$ws.Range("H85").Value = 1100.5
$ws.Range("J85").Value = 134
$ws.Range("L85").Value = 134
$ws.Range("N85").Value = -2630
$ws.Range("H113").Value = 2841
$ws.Range("I113").Value = 2481.5
$ws.Range("K113").Value = 2481.5
$ws.Range("M113").Value = -311.5
$ws.Range("H122").Value = 1957.2858
$ws.Range("I122").Value = 1825.5
$ws.Range("J122").Value = 2133
$ws.Range("K122").Value = 5476.5
$ws.Range("L122").Value = 6399
$ws.Range("M122").Value = -3026.5
$ws.Range("N122").Value = -11299
$ws.Range("H132").Value = 5599.2
$ws.Range("I132").Value = 5004
$ws.Range("J132").Value = 5748
$ws.Range("K132").Value = 15012
$ws.Range("L132").Value = 17244
$ws.Range("M132").Value = -12482
$ws.Range("N132").Value = -22304
$ws.Range("H136").Value = 19310
$ws.Range("I136").Value = 13172.2
$ws.Range("K136").Value = 39516.60000000001
$ws.Range("M136").Value = -36966.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1623.909
$ws.Range("I81").Value = 1885.8889
$ws.Range("J81").Value = 445
$ws.Range("K81").Value = 3771.7778
$ws.Range("L81").Value = 890
$ws.Range("M81").Value = -2710.7778
$ws.Range("N81").Value = -3012
$ws.Range("H84").Value = 1623.909
$ws.Range("I84").Value = 1885.8889
$ws.Range("J84").Value = 445
$ws.Range("K84").Value = 18858.889
$ws.Range("L84").Value = 4450
$ws.Range("M84").Value = -13554.889
$ws.Range("N84").Value = -15058
$ws.Range("H132").Value = 2524.889
$ws.Range("I132").Value = 1850.6428
$ws.Range("K132").Value = 5551.928400000001
$ws.Range("M132").Value = -3021.928400000001
$ws.Range("H136").Value = 14712.111
$ws.Range("I136").Value = 14712.111
$ws.Range("K136").Value = 44136.333
$ws.Range("M136").Value = -41586.333
